$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a new attendance row's date. A plain string Value
# assignment of "2024-04-13" would be auto-parsed by COM into a date
# serial number, but the target keeps it as literal text, so force the
# cell to Text format first, then restore the default "Normal" style so
# no extra formatting sticks to the cell (matches the source row, which
# carries no explicit style on the new row's cells).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2024-04-13"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "20:59:53"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "shugi"
$ws.Range("E2").Value = "21:00:00"
$ws.Range("F2").Value = "00:00:07"
